$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1) Insert a new row at 52 (pushes old rows 52-57 -> 53-58, i.e. the
#    signature block moves from rows 56/57 to rows 57/58), to make room
#    for one extra "Periodo Mora" entry in the data table.
# ------------------------------------------------------------------
$ws.Rows("52:52").Insert()

# New row 52 should inherit the "closing" (last-row) border formatting
# that row 51 currently has.
$ws.Range("B51:J51").Copy()
$ws.Range("B52:J52").PasteSpecial(-4122)  # xlPasteFormats

# Row 51 is no longer the last row of the table, so it goes back to the
# regular data-row formatting (same as row 50).
$ws.Range("B50:J50").Copy()
$ws.Range("B51:J51").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2) Fill in B52:D52 / F52:G52 the same way every other data row is
#    populated (same worker, same mora/salary values).
# ------------------------------------------------------------------
$ws.Range("B52").Value = "CC"
$ws.Range("C52").Value = "9293613"
$ws.Range("D52").Value = "JOSIAS DE JESUS HURTADO TORRES"
$ws.Range("F52").Value = 40000
$ws.Range("G52").Value = 1000000

# ------------------------------------------------------------------
# 3) Re-populate B16:D51 (worker data repeated down the new row too)
#    and the "Periodo Mora" column E16:E52 with the updated list,
#    newest period first (a new period, 2507, was added).
# ------------------------------------------------------------------
$periods = @(
    "2507","2506","2505","2504","2503","2502","2501",
    "2412","2411","2410","2409","2408","2407","2406","2405","2404","2403","2402","2401",
    "2312","2311","2310","2309","2308","2307","2306","2305","2304","2303","2302","2301",
    "2212","2211","2210","2209","2208","2207"
)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 16 + $i
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = "9293613"
    $ws.Range("D$r").Value = "JOSIAS DE JESUS HURTADO TORRES"
    $ws.Range("E$r").Value = $periods[$i]
    $ws.Range("F$r").Value = 40000
    $ws.Range("G$r").Value = 1000000
}

# ------------------------------------------------------------------
# 4) Update the summary figures: one more period => Cant. Periodos
#    36 -> 37, and total mora grows by one more 40000 period.
# ------------------------------------------------------------------
$ws.Range("F13").Value = 37
$ws.Range("E11").Value = 1480000

Write-Host "Edit complete"
